# Update the Drought Action Start Month (and the linked Decision start
# month) on the "UI" sheet from June (6) to July (7).
$wb = $excel.ActiveWorkbook

$uiSheet = $wb.Worksheets.Item("UI")

# F9  = "Start Month (Drought)"
# F10 = "Start Month (Decision)"
$uiSheet.Range("F9").Value = 7
$uiSheet.Range("F10").Value = 7

# Recalculate so every dependent formula's cached value is refreshed.
$excel.Calculate()

# Refresh every chart so its cached numCache snapshot picks up the new
# dependent values too.
foreach ($sheet in $wb.Worksheets) {
    foreach ($co in $sheet.ChartObjects()) {
        $co.Chart.Refresh()
    }
}

# Make "UI" the active sheet/tab, with F10 selected (matches the author's
# final on-screen state after making the edit).
$uiSheet.Activate()
$uiSheet.Range("F10").Select()
